$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The capsule name changed from "Challengers" to "Contenders" everywhere it
# is used (A1 and A2 share this string).
$ws.Range("A1:A2").Value = "Paris 2023 Contenders Sticker Capsule"

# A3 and A4 no longer carry that text - clear their contents (rows remain,
# cells become blank).
$ws.Range("A3:A4").ClearContents()
